$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(32, 8).Value = 277.1111
$ws.Cells.Item(32, 9).Value = 287.5
$ws.Cells.Item(32, 10).Value = 274.14285
$ws.Cells.Item(32, 11).Value = 287.5
$ws.Cells.Item(32, 12).Value = 274.14285
$ws.Cells.Item(32, 13).Value = 38.5
$ws.Cells.Item(32, 14).Value = -926.14285

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(33, 8).Value = 300194.06
$ws.Cells.Item(33, 9).Value = 433.66666
$ws.Cells.Item(33, 11).Value = 433.66666
$ws.Cells.Item(33, 13).Value = -204.66666

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(62, 8).Value = 2996
$ws.Cells.Item(62, 9).Value = 2500
$ws.Cells.Item(62, 10).Value = 4980
$ws.Cells.Item(62, 11).Value = 2500
$ws.Cells.Item(62, 12).Value = 4980
$ws.Cells.Item(62, 13).Value = -1876
$ws.Cells.Item(62, 14).Value = -6228

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(65, 8).Value = 2996
$ws.Cells.Item(65, 9).Value = 2500
$ws.Cells.Item(65, 10).Value = 4980
$ws.Cells.Item(65, 11).Value = 12500
$ws.Cells.Item(65, 12).Value = 24900
$ws.Cells.Item(65, 13).Value = -9380
$ws.Cells.Item(65, 14).Value = -31140

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(70, 8).Value = 1099.5312
$ws.Cells.Item(70, 9).Value = 1220.7142
$ws.Cells.Item(70, 11).Value = 3662.1426
$ws.Cells.Item(70, 13).Value = -3392.1426

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(73, 8).Value = 1099.5312
$ws.Cells.Item(73, 9).Value = 1220.7142
$ws.Cells.Item(73, 11).Value = 3662.1426
$ws.Cells.Item(73, 13).Value = -2726.1426

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(86, 8).Value = 9339
$ws.Cells.Item(86, 9).Value = 2566.6667
$ws.Cells.Item(86, 10).Value = 19497.5
$ws.Cells.Item(86, 11).Value = 2566.6667
$ws.Cells.Item(86, 12).Value = 19497.5
$ws.Cells.Item(86, 13).Value = -1443.6667
$ws.Cells.Item(86, 14).Value = -21743.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(89, 8).Value = 9339
$ws.Cells.Item(89, 9).Value = 2566.6667
$ws.Cells.Item(89, 10).Value = 19497.5
$ws.Cells.Item(89, 11).Value = 12833.3335
$ws.Cells.Item(89, 12).Value = 97487.5
$ws.Cells.Item(89, 13).Value = -7217.333500000001
$ws.Cells.Item(89, 14).Value = -108719.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(98, 8).Value = 1159.25
$ws.Cells.Item(98, 9).Value = 1283.7
$ws.Cells.Item(98, 10).Value = 537
$ws.Cells.Item(98, 11).Value = 1283.7
$ws.Cells.Item(98, 12).Value = 537
$ws.Cells.Item(98, 13).Value = 214.3
$ws.Cells.Item(98, 14).Value = -3533

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(113, 8).Value = 102141.1
$ws.Cells.Item(113, 9).Value = 202281
$ws.Cells.Item(113, 10).Value = 2001.2
$ws.Cells.Item(113, 11).Value = 202281
$ws.Cells.Item(113, 12).Value = 2001.2
$ws.Cells.Item(113, 13).Value = -199027
$ws.Cells.Item(113, 14).Value = -8509.200000000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(122, 8).Value = 1159.25
$ws.Cells.Item(122, 9).Value = 1283.7
$ws.Cells.Item(122, 10).Value = 537
$ws.Cells.Item(122, 11).Value = 3851.1
$ws.Cells.Item(122, 12).Value = 1611
$ws.Cells.Item(122, 13).Value = -1401.1
$ws.Cells.Item(122, 14).Value = -6511

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(129, 8).Value = 2150.5789
$ws.Cells.Item(129, 9).Value = 5485.55
$ws.Cells.Item(129, 10).Value = 959.5179000000001
$ws.Cells.Item(129, 11).Value = 16456.65
$ws.Cells.Item(129, 12).Value = 2878.5537
$ws.Cells.Item(129, 13).Value = -11456.65
$ws.Cells.Item(129, 14).Value = -12878.5537

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(138, 8).Value = 1666.6538
$ws.Cells.Item(138, 9).Value = 1446.1072
$ws.Cells.Item(138, 10).Value = 1923.9584
$ws.Cells.Item(138, 11).Value = 4338.321599999999
$ws.Cells.Item(138, 12).Value = 5771.8752
$ws.Cells.Item(138, 13).Value = 801.6784000000007
$ws.Cells.Item(138, 14).Value = -16051.8752

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 21500.988
$ws.Cells.Item(32, 9).Value = 4799.711
$ws.Cells.Item(32, 11).Value = 4799.711
$ws.Cells.Item(32, 13).Value = -4512.711

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(132, 8).Value = 1799.5319
$ws.Cells.Item(132, 9).Value = 1536.122
$ws.Cells.Item(132, 10).Value = 3599.5
$ws.Cells.Item(132, 11).Value = 4608.366
$ws.Cells.Item(132, 12).Value = 10798.5
$ws.Cells.Item(132, 13).Value = -2078.366
$ws.Cells.Item(132, 14).Value = -15858.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(139, 8).Value = 30373.572
$ws.Cells.Item(139, 10).Value = 32102.5
$ws.Cells.Item(139, 12).Value = 32102.5
$ws.Cells.Item(139, 14).Value = -42382.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(59, 8).Value = 67990
$ws.Cells.Item(59, 10).Value = 67990
$ws.Cells.Item(59, 12).Value = 67990
$ws.Cells.Item(59, 14).Value = -69684

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(137, 8).Value = 59849.332
$ws.Cells.Item(137, 10).Value = 59849.332
$ws.Cells.Item(137, 12).Value = 59849.332
$ws.Cells.Item(137, 14).Value = -70049.33199999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(140, 8).Value = 66995
$ws.Cells.Item(140, 10).Value = 66995
$ws.Cells.Item(140, 12).Value = 66995
$ws.Cells.Item(140, 14).Value = -77355

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(138, 8).Value = 56570.715
$ws.Cells.Item(138, 10).Value = 56570.715
$ws.Cells.Item(138, 12).Value = 56570.715
$ws.Cells.Item(138, 14).Value = -66850.715

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(131, 8).Value = 799.0700000000001
$ws.Cells.Item(131, 9).Value = 292.46155
$ws.Cells.Item(131, 10).Value = 874.77014
$ws.Cells.Item(131, 11).Value = 877.38465
$ws.Cells.Item(131, 12).Value = 2624.31042
$ws.Cells.Item(131, 13).Value = 4162.61535
$ws.Cells.Item(131, 14).Value = -12704.31042

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(141, 8).Value = 11990
$ws.Cells.Item(141, 9).Value = 21098
$ws.Cells.Item(141, 10).Value = 4400
$ws.Cells.Item(141, 11).Value = 63294
$ws.Cells.Item(141, 12).Value = 13200
$ws.Cells.Item(141, 13).Value = -58114
$ws.Cells.Item(141, 14).Value = -23560

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(64, 8).Value = 44695.332
$ws.Cells.Item(64, 10).Value = 44695.332
$ws.Cells.Item(64, 12).Value = 44695.332
$ws.Cells.Item(64, 14).Value = -45191.332

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(67, 8).Value = 44695.332
$ws.Cells.Item(67, 10).Value = 44695.332
$ws.Cells.Item(67, 12).Value = 44695.332
$ws.Cells.Item(67, 14).Value = -46411.332

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 61195.43
$ws.Cells.Item(70, 9).Value = 87287.875
$ws.Cells.Item(70, 10).Value = 4266.4546
$ws.Cells.Item(70, 11).Value = 87287.875
$ws.Cells.Item(70, 12).Value = 4266.4546
$ws.Cells.Item(70, 13).Value = -87017.875
$ws.Cells.Item(70, 14).Value = -4806.4546

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(73, 8).Value = 61195.43
$ws.Cells.Item(73, 9).Value = 87287.875
$ws.Cells.Item(73, 10).Value = 4266.4546
$ws.Cells.Item(73, 11).Value = 87287.875
$ws.Cells.Item(73, 12).Value = 4266.4546
$ws.Cells.Item(73, 13).Value = -86351.875
$ws.Cells.Item(73, 14).Value = -6138.4546

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 2275.6743
$ws.Cells.Item(132, 9).Value = 1894.3143
$ws.Cells.Item(132, 10).Value = 3944.125
$ws.Cells.Item(132, 11).Value = 5682.9429
$ws.Cells.Item(132, 12).Value = 11832.375
$ws.Cells.Item(132, 13).Value = -3152.9429
$ws.Cells.Item(132, 14).Value = -16892.375

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(139, 8).Value = 58611.2
$ws.Cells.Item(139, 10).Value = 58611.2
$ws.Cells.Item(139, 12).Value = 58611.2
$ws.Cells.Item(139, 14).Value = -68891.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(140, 8).Value = 97859.91
$ws.Cells.Item(140, 10).Value = 97859.91
$ws.Cells.Item(140, 12).Value = 97859.91
$ws.Cells.Item(140, 14).Value = -108219.91

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(132, 8).Value = 4461.6
$ws.Cells.Item(132, 9).Value = 4519.1665
$ws.Cells.Item(132, 10).Value = 3943.5
$ws.Cells.Item(132, 11).Value = 13557.4995
$ws.Cells.Item(132, 12).Value = 11830.5
$ws.Cells.Item(132, 13).Value = -11027.4995
$ws.Cells.Item(132, 14).Value = -16890.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 200569.8
$ws.Cells.Item(81, 9).Value = 125599.875
$ws.Cells.Item(81, 10).Value = 500449.5
$ws.Cells.Item(81, 11).Value = 251199.75
$ws.Cells.Item(81, 12).Value = 1000899
$ws.Cells.Item(81, 13).Value = -250138.75
$ws.Cells.Item(81, 14).Value = -1003021

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(84, 8).Value = 200569.8
$ws.Cells.Item(84, 9).Value = 125599.875
$ws.Cells.Item(84, 10).Value = 500449.5
$ws.Cells.Item(84, 11).Value = 1255998.75
$ws.Cells.Item(84, 12).Value = 5004495
$ws.Cells.Item(84, 13).Value = -1250694.75
$ws.Cells.Item(84, 14).Value = -5015103

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(136, 8).Value = 524.0833
$ws.Cells.Item(136, 9).Value = 377.39285
$ws.Cells.Item(136, 10).Value = 1037.5
$ws.Cells.Item(136, 11).Value = 1132.17855
$ws.Cells.Item(136, 12).Value = 3112.5
$ws.Cells.Item(136, 13).Value = 1417.82145
$ws.Cells.Item(136, 14).Value = -8212.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(138, 8).Value = 66132.5
$ws.Cells.Item(138, 10).Value = 66132.5
$ws.Cells.Item(138, 12).Value = 66132.5
$ws.Cells.Item(138, 14).Value = -76412.5
